$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a Text number format ("@", numFmtId 49) to the whole used range so
# that phone numbers read from Excel keep their leading "+" / formatting
# instead of being coerced into numbers.
$ws.Range("A1:E3").NumberFormat = "@"

# Correct the phone number values (previously misread) to their correct,
# textual, "+"-prefixed representations.
$ws.Range("D3").Value = "+17324061005"
$ws.Range("D2").Value = "+11231231234"

# Update the active selection to match the author's saved state.
$ws.Range("E2").Select()
